$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the shared string used by C4 from "smart_attribute" to "smart_attribute_state"
$ws.Range("C4").Value = "smart_attribute_state"

# 2. Move the active cell / selection from C3 to C4
$ws.Range("C4").Select()

# 3. Nudge the widths of columns B, C and D slightly wider (closest values the
#    host's 1/6-character quantized ColumnWidth can represent).
$ws.Columns.Item(2).ColumnWidth = 15.5
$ws.Columns.Item(3).ColumnWidth = 21.6667
$ws.Columns.Item(4).ColumnWidth = 21.6667
